$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-7: swap in the new "(SYS_USER3)" style API paths
#     and refresh the query-string / validation text for rows 6-7 ---
$ws.Range("D2").Value = "/entitlements/(SYS_USER3)"
$ws.Range("D3").Value = "/entitlements/(SYS_USER3)1"

$ws.Range("D6").Value = "/entitlements/(SYS_USER3)/entitled"
$ws.Range("G6").Value = "?entitilement=214504 OR entitilement=215802"
$ws.Range("J6").Value = "status=403||errorCode=403.1.1||errorMessage=User does not have sku"

$ws.Range("D7").Value = "/entitlements/(SYS_USER3)1/entitled"
$ws.Range("G7").Value = "?entitilement=214504 OR entitilement=215802"

# --- New test case rows 8, 9 and 10 ---
$ws.Rows.Item(8).RowHeight = 45
$ws.Range("A8").Value = "OPQA-3852"
$ws.Range("B8").Value = "Verify that to get specific entitlement of the user by passing valid truid and valid entitlement name"
$ws.Range("C8").Value = "1PENTITLEMENTS"
$ws.Range("D8").Value = "/entitlements/filter/(SYS_USER3)/DRA_TARGET_DRUG"
$ws.Range("E8").Value = "GET"
$ws.Range("J8").Value = "status=200||skus=DRA_TARGET_DRUG||X-1P-ENT=DRA"
$ws.Range("L8").Clear()

$ws.Rows.Item(9).RowHeight = 45
$ws.Range("A9").Value = "OPQA-3853"
$ws.Range("B9").Value = "Verify that to get morethan one specific entitlement of the user by passing valid truid and valid entitlement names saperated by comma(,)"
$ws.Range("C9").Value = "1PENTITLEMENTS"
$ws.Range("D9").Value = "/entitlements/filter/(SYS_USER3)/DRA_TARGET_DRUG,IPA_TEST_SKU"
$ws.Range("E9").Value = "GET"
$ws.Range("J9").Value = "status=200||skus=DRA_TARGET_DRUG||skus=IPA_TEST_SKU||X-1P-ENT=DRA"
$ws.Range("L9").Clear()

$ws.Rows.Item(10).RowHeight = 45
$ws.Range("A10").Value = "OPQA-3856"
$ws.Range("B10").Value = "Verify that to get error status by passing invalid truid and valid entitlement id"
$ws.Range("C10").Value = "1PENTITLEMENTS"
$ws.Range("D10").Value = "/entitlements/filter/(SYS_USER3)1/DRA_TARGET_DRUG"
$ws.Range("E10").Value = "GET"
$ws.Range("J10").Value = "status=403||errorCode=403.1.8||errorMessage=Unknowen user"
$ws.Range("L10").Clear()

# --- Column D widened to fit the new, longer API path strings ---
$ws.Columns.Item(4).ColumnWidth = 74.5

# --- Selection now spans the newly added rows ---
$excel.Goto($ws.Range("L2:L10"))
